$wb = $excel.ActiveWorkbook

# Activate the "6x6" sheet (it was already the tab-selected sheet) and
# update its view: scroll back to the top-left (clears topLeftCell) and
# change the selection from the whole-column selection A1:XFD1048576 down
# to the first 8 rows, A1:XFD8.
$ws = $wb.Worksheets.Item("6x6")
$ws.Activate() | Out-Null
$ws.Rows("1:8").Select() | Out-Null

# Force a full recalculation so the volatile RAND()-based cached values on
# the "rand" sheet are refreshed (matches the workbook having been
# recalculated/resaved).
$excel.CalculateFull() | Out-Null
